$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.129.65"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "2.549.00"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "516.23"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "141.08"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").Value = "2.553.83"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "2.998.98"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").Value = "57.141.40"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "20.00"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").Value = "2.570.85"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "330.92"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "6.15"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "64.73"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "0.398"
$ws.Range("E27").Value = "  -4.98%  "
$ws.Range("D28").Value = "2.654.85"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  -8.18%  "
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -6.28%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.46"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "148.42"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("D38").Value = "0.832"
$ws.Range("E38").Value = "  -7.56%  "
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("D40").Value = "0.821"
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").Value = "10.62"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "0.0949"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "265.06"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").Value = "0.578"
$ws.Range("E47").Value = "  -6.15%  "
$ws.Range("D48").Value = "18.56"
$ws.Range("E48").Value = "  -6.26%  "
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("D50").Value = "1.959.87"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("E51").Value = "  -4.53%  "
